$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename the Average Rainfall / Average Temp / Wind Max columns
$ws.Range("C1").Value = "Average_Rainfall "
$ws.Range("D1").Value = "Average_Temp"
$ws.Range("E1").Value = "Wind_Max "

# Abbreviate month names in column B (rows 2-13), May stays "May "
$ws.Range("B2").Value = "Jan"
$ws.Range("B3").Value = "Feb"
$ws.Range("B4").Value = "Mar"
$ws.Range("B5").Value = "Apr"
$ws.Range("B6").Value = "May "
$ws.Range("B7").Value = "Jun"
$ws.Range("B8").Value = "Jul"
$ws.Range("B9").Value = "Aug"
$ws.Range("B10").Value = "Sep"
$ws.Range("B11").Value = "Oct"
$ws.Range("B12").Value = "Nov "
$ws.Range("B13").Value = "Dec "

# Update selection to match the target state
$ws.Range("D7").Select()
